# ajout cyclo tour alsace
# Replace the old row 40 (reported/postponed Seppois entry, now obsolete)
# with the new "Etape Cyclo du Tour Alsace" event, and clear the
# now-unused last column (report note) for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A40").Value = "Sam 24 Juillet"
$ws.Range("B40").Value = "Etape Cyclo du Tour Alsace"
$ws.Range("C40").Value = "FSGT"
$ws.Range("D40").Value = "Randonnée"
$ws.Range("E40").Value = "https://www.touralsace.fr/boutique/etape-cyclo/"
$ws.Range("F40").ClearContents()
$ws.Range("G40").ClearContents()

# Update the view: scroll position / current selection as saved in the
# workbook when the edit was made.
$ws.Range("E41").Select()
